$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header change
$ws.Range("C1").Value = "runtime"

# Row 2
$ws.Range("C2").Value = "Node.js"
$ws.Range("D2").Value = 58729
$ws.Range("E2").Value = 37.48
$ws.Range("F2").Value = 11.19
$ws.Range("G2").Value = 24.71
$ws.Range("H2").Value = 34.93
$ws.Range("I2").Value = 67.53
$ws.Range("J2").Value = 544.77

# Row 3
$ws.Range("C3").Value = "Python"
$ws.Range("D3").Value = 58698
$ws.Range("E3").Value = 35.86
$ws.Range("F3").Value = 8.380000000000001
$ws.Range("G3").Value = 23.81
$ws.Range("H3").Value = 33.94
$ws.Range("I3").Value = 61.2
$ws.Range("J3").Value = 1040.8

# Row 4
$ws.Range("C4").Value = "Node.js"
$ws.Range("D4").Value = 56696
$ws.Range("E4").Value = 19.81
$ws.Range("F4").Value = 5.84
$ws.Range("G4").Value = 14.43
$ws.Range("H4").Value = 18.21
$ws.Range("I4").Value = 40.3
$ws.Range("J4").Value = 255.48

# Row 5
$ws.Range("C5").Value = "Python"
$ws.Range("D5").Value = 56577
$ws.Range("E5").Value = 21.45
$ws.Range("F5").Value = 6.1
$ws.Range("G5").Value = 15.32
$ws.Range("H5").Value = 19.86
$ws.Range("I5").Value = 43.46
$ws.Range("J5").Value = 190.91

# Row 6
$ws.Range("C6").Value = "Golang"
$ws.Range("D6").Value = 57408
$ws.Range("E6").Value = 151.82
$ws.Range("F6").Value = 130.37
$ws.Range("G6").Value = 37.34
$ws.Range("H6").Value = 191.41
$ws.Range("I6").Value = 424.88
$ws.Range("J6").Value = 1128.18

# Row 7
$ws.Range("C7").Value = "Node.js"
$ws.Range("D7").Value = 57384
$ws.Range("E7").Value = 152.44
$ws.Range("F7").Value = 131.36
$ws.Range("G7").Value = 39.37
$ws.Range("H7").Value = 191.94
$ws.Range("I7").Value = 425.42
$ws.Range("J7").Value = 5102.82

# Row 8
$ws.Range("C8").Value = "Golang"
$ws.Range("D8").Value = 58777
$ws.Range("E8").Value = 42.45
$ws.Range("F8").Value = 7.44
$ws.Range("G8").Value = 37.46
$ws.Range("H8").Value = 40.72
$ws.Range("I8").Value = 72.7
$ws.Range("J8").Value = 179.21

# Row 9
$ws.Range("C9").Value = "Node.js"
$ws.Range("D9").Value = 58772
$ws.Range("E9").Value = 43.28
$ws.Range("F9").Value = 19.59
$ws.Range("G9").Value = 37.41
$ws.Range("H9").Value = 41.09
$ws.Range("I9").Value = 74.44
$ws.Range("J9").Value = 3549.11

# Row 10
$ws.Range("C10").Value = "Python"
$ws.Range("D10").Value = 58746
$ws.Range("E10").Value = 43.11
$ws.Range("F10").Value = 7.48
$ws.Range("G10").Value = 38.41
$ws.Range("H10").Value = 41.4
$ws.Range("I10").Value = 72.75
$ws.Range("J10").Value = 334.4
